$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 19.30324074074072
$ws.Range("I5").Value = 19.30324074074072
$ws.Range("I8").Value = 12.67039049919483
$ws.Range("I9").Value = 12.67039049919483
$ws.Range("I10").Value = 13.17361111111111
$ws.Range("I11").Value = 13.17361111111111
$ws.Range("I12").Value = 12.92654320987656
$ws.Range("I13").Value = 12.92654320987656
$ws.Range("I14").Value = -1.226851851851833
$ws.Range("I15").Value = -1.226851851851833
$ws.Range("I16").Value = 13.0158303464755
$ws.Range("I17").Value = 13.0158303464755
$ws.Range("I20").Value = 1.925925925925943
$ws.Range("I21").Value = 1.925925925925943
$ws.Range("I22").Value = -1.819444444444444
$ws.Range("I23").Value = -1.819444444444444
$ws.Range("I24").Value = 13.0158303464755
$ws.Range("I25").Value = 13.0158303464755
$ws.Range("I26").Value = 5.486111111111112
$ws.Range("I27").Value = 5.486111111111112
$ws.Range("I30").Value = -1.226851851851833
$ws.Range("I31").Value = -1.226851851851833
$ws.Range("I32").Value = 3.38888888888889
$ws.Range("I33").Value = 3.38888888888889
$ws.Range("I38").Value = 16.86342592592595
$ws.Range("I39").Value = 16.86342592592595

$ws.Range("N4").Value = 1.352319749654237
$ws.Range("O4").Value = 1.423118895050623
$ws.Range("N8").Value = 1.300269876134972
$ws.Range("O8").Value = 1.364969709970879
$ws.Range("N10").Value = 1.304077921028169
$ws.Range("O10").Value = 1.369214264257821
$ws.Range("N12").Value = 1.302205489329493
$ws.Range("O12").Value = 1.367127007643996
$ws.Range("N14").Value = 1.203236793039155
$ws.Range("O14").Value = 1.257328254301852
$ws.Range("N16").Value = 1.302881541082627
$ws.Range("O16").Value = 1.367880580392128
$ws.Range("N20").Value = 1.223958122597613
$ws.Range("O20").Value = 1.280232184891932
$ws.Range("N22").Value = 1.199420117463385
$ws.Range("O22").Value = 1.253114445055376
$ws.Range("N24").Value = 1.302881541082627
$ws.Range("O24").Value = 1.367880580392128
$ws.Range("N26").Value = 1.248232108317215
$ws.Range("O26").Value = 1.307120032773454
$ws.Range("N30").Value = 1.203236793039155
$ws.Range("O30").Value = 1.257328254301852
$ws.Range("N32").Value = 1.233817681248088
$ws.Range("O32").Value = 1.291146001942376
$ws.Range("N38").Value = 1.332696358504853
$ws.Range("O38").Value = 1.401162263046183
